$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.674985639065483
$ws.Cells.Item(2, 4).Value = 0.1739349866106039
$ws.Cells.Item(2, 5).Value = 0.253121296826734
$ws.Cells.Item(2, 6).Value = 2.05332449930593
$ws.Cells.Item(2, 7).Value = 1.48442475116633
$ws.Cells.Item(2, 8).Value = 1.326505145020036
$ws.Cells.Item(2, 10).Value = 0.3827937704926114
$ws.Cells.Item(2, 12).Value = 0.5934112428884362
$ws.Cells.Item(2, 13).Value = 0.4907106959707761

$ws.Cells.Item(3, 2).Value = 1.596616672991161
$ws.Cells.Item(3, 4).Value = 0.1671798858679381
$ws.Cells.Item(3, 5).Value = 0.2419745318911026
$ws.Cells.Item(3, 6).Value = 2.068119898121935
$ws.Cells.Item(3, 7).Value = 1.479313313549753
$ws.Cells.Item(3, 8).Value = 1.332673545421528
$ws.Cells.Item(3, 10).Value = 0.3647985006520145
$ws.Cells.Item(3, 12).Value = 0.5378407109201646
$ws.Cells.Item(3, 13).Value = 0.4579375082744335

$ws.Cells.Item(4, 2).Value = 1.54893231772823
$ws.Cells.Item(4, 4).Value = 0.162996975481704
$ws.Cells.Item(4, 5).Value = 0.2351104525886072
$ws.Cells.Item(4, 6).Value = 2.07889269964349
$ws.Cells.Item(4, 7).Value = 1.47762503735737
$ws.Cells.Item(4, 8).Value = 1.337446048174797
$ws.Cells.Item(4, 10).Value = 0.3537419335614089
$ws.Cells.Item(4, 12).Value = 0.5037233380477346
$ws.Cells.Item(4, 13).Value = 0.4378910731713646

$ws.Cells.Item(5, 2).Value = 1.529610490772058
$ws.Cells.Item(5, 4).Value = 0.1612834959818912
$ws.Cells.Item(5, 5).Value = 0.2323083969131332
$ws.Cells.Item(5, 6).Value = 2.083705882097973
$ws.Cells.Item(5, 7).Value = 1.477299788878682
$ws.Cells.Item(5, 8).Value = 1.339637871516047
$ws.Cells.Item(5, 10).Value = 0.3492347173716297
$ws.Cells.Item(5, 12).Value = 0.4898214643109782
$ws.Cells.Item(5, 13).Value = 0.4297414664282613

$ws.Cells.Item(6, 2).Value = 1.526408776418208
$ws.Cells.Item(6, 4).Value = 0.1609984346206517
$ws.Cells.Item(6, 5).Value = 0.2318428261747627
$ws.Cells.Item(6, 6).Value = 2.084530626116688
$ws.Cells.Item(6, 7).Value = 1.47726763403081
$ws.Cells.Item(6, 8).Value = 1.340016718771651
$ws.Cells.Item(6, 10).Value = 0.3484862109732632
$ws.Cells.Item(6, 12).Value = 0.487513156595071
$ws.Cells.Item(6, 13).Value = 0.4283894147474001

$ws.Cells.Item(7, 2).Value = 1.548671290786785
$ws.Cells.Item(7, 4).Value = 0.1629739030281598
$ws.Cells.Item(7, 5).Value = 0.2350726827086547
$ws.Cells.Item(7, 6).Value = 2.07895590031363
$ws.Cells.Item(7, 7).Value = 1.477619184781588
$ws.Cells.Item(7, 8).Value = 1.337474608783381
$ws.Cells.Item(7, 10).Value = 0.3536811536589681
$ws.Cells.Item(7, 12).Value = 0.503535846925331
$ws.Cells.Item(7, 13).Value = 0.4377810855077584

$ws.Cells.Item(8, 2).Value = 1.647874207200061
$ws.Cells.Item(8, 4).Value = 0.1716131145815893
$ws.Cells.Item(8, 5).Value = 0.2492820949962464
$ws.Cells.Item(8, 6).Value = 2.058074701318546
$ws.Cells.Item(8, 7).Value = 1.482360182996558
$ws.Cells.Item(8, 8).Value = 1.32842712157634
$ws.Cells.Item(8, 10).Value = 0.376590699896596
$ws.Cells.Item(8, 12).Value = 0.5742501072764696
$ws.Cells.Item(8, 13).Value = 0.4793947775182019

$ws.Cells.Item(9, 2).Value = 1.845841062963927
$ws.Cells.Item(9, 4).Value = 0.1882777558568733
$ws.Cells.Item(9, 5).Value = 0.2769848635144001
$ws.Cells.Item(9, 6).Value = 2.030586713591006
$ws.Cells.Item(9, 7).Value = 1.503253315598641
$ws.Cells.Item(9, 8).Value = 1.318534853009425
$ws.Cells.Item(9, 10).Value = 0.4214476361764525
$ws.Cells.Item(9, 12).Value = 0.7129343630377321
$ws.Cells.Item(9, 13).Value = 0.561599114294367

$ws.Cells.Item(10, 2).Value = 1.993373260071621
$ws.Cells.Item(10, 4).Value = 0.2003576876318931
$ws.Cells.Item(10, 5).Value = 0.2972366319412387
$ws.Cells.Item(10, 6).Value = 2.018687056306121
$ws.Cells.Item(10, 7).Value = 1.525801320448721
$ws.Cells.Item(10, 8).Value = 1.316101015179214
$ws.Cells.Item(10, 10).Value = 0.454352737727433
$ws.Cells.Item(10, 12).Value = 0.8148308189254294
$ws.Cells.Item(10, 13).Value = 0.6223589608455029

$ws.Cells.Item(11, 2).Value = 2.060942957234545
$ws.Cells.Item(11, 4).Value = 0.2058188602901794
$ws.Cells.Item(11, 5).Value = 0.3064272829191381
$ws.Cells.Item(11, 6).Value = 2.015094631159286
$ws.Cells.Item(11, 7).Value = 1.537650019212549
$ws.Cells.Item(11, 8).Value = 1.316054132879827
$ws.Cells.Item(11, 10).Value = 0.4693093601215992
$ws.Cells.Item(11, 12).Value = 0.8611876832263476
$ws.Cells.Item(11, 13).Value = 0.6500796499821462

$ws.Cells.Item(12, 2).Value = 2.086595218679463
$ws.Cells.Item(12, 4).Value = 0.2078820450846308
$ws.Cells.Item(12, 5).Value = 0.3099043121678307
$ws.Cells.Item(12, 6).Value = 2.013997685810409
$ws.Cells.Item(12, 7).Value = 1.542367880732542
$ws.Cells.Item(12, 8).Value = 1.316189696380548
$ws.Cells.Item(12, 10).Value = 0.4749711100330671
$ws.Cells.Item(12, 12).Value = 0.878742175002003
$ws.Cells.Item(12, 13).Value = 0.6605882642996619

$ws.Cells.Item(13, 2).Value = 2.081067659314158
$ws.Cells.Item(13, 4).Value = 0.2074379158325996
$ws.Cells.Item(13, 5).Value = 0.3091556188788616
$ws.Cells.Item(13, 6).Value = 2.014222190984867
$ws.Cells.Item(13, 7).Value = 1.541341495777544
$ws.Cells.Item(13, 8).Value = 1.316153667944064
$ws.Cells.Item(13, 10).Value = 0.4737518438446671
$ws.Cells.Item(13, 12).Value = 0.8749614999666164
$ws.Cells.Item(13, 13).Value = 0.6583245428950306

$ws.Cells.Item(14, 2).Value = 2.063052082999548
$ws.Cells.Item(14, 4).Value = 0.2059886964861875
$ws.Cells.Item(14, 5).Value = 0.3067134063729426
$ws.Cells.Item(14, 6).Value = 2.014999096111751
$ws.Cells.Item(14, 7).Value = 1.538033518184506
$ws.Cells.Item(14, 8).Value = 1.316062208184178
$ws.Cells.Item(14, 10).Value = 0.4697751972139201
$ws.Cells.Item(14, 12).Value = 0.8626319008049563
$ws.Cells.Item(14, 13).Value = 0.6509439729157407

$ws.Cells.Item(15, 2).Value = 2.052025481762712
$ws.Cells.Item(15, 4).Value = 0.2051003787459109
$ws.Cells.Item(15, 5).Value = 0.3052170515559212
$ws.Cells.Item(15, 6).Value = 2.015509326718046
$ws.Cells.Item(15, 7).Value = 1.536037435896048
$ws.Cells.Item(15, 8).Value = 1.316026178140703
$ws.Cells.Item(15, 10).Value = 0.4673391184558682
$ws.Cells.Item(15, 12).Value = 0.8550796746312699
$ws.Cells.Item(15, 13).Value = 0.6464246378922667

$ws.Cells.Item(16, 2).Value = 1.988966559235678
$ws.Cells.Item(16, 4).Value = 0.2000001065175212
$ws.Cells.Item(16, 5).Value = 0.2966355494657051
$ws.Cells.Item(16, 6).Value = 2.018958605296945
$ws.Cells.Item(16, 7).Value = 1.525059189578485
$ws.Cells.Item(16, 8).Value = 1.316125487102795
$ws.Cells.Item(16, 10).Value = 0.4533750217022003
$ws.Cells.Item(16, 12).Value = 0.8118013300168059
$ws.Cells.Item(16, 13).Value = 0.6205489541804781

$ws.Cells.Item(17, 2).Value = 1.950398473825715
$ws.Cells.Item(17, 4).Value = 0.1968625712443384
$ws.Cells.Item(17, 5).Value = 0.2913653663647437
$ws.Cells.Item(17, 6).Value = 2.021542109752062
$ws.Cells.Item(17, 7).Value = 1.518733545647649
$ws.Cells.Item(17, 8).Value = 1.316458575999036
$ws.Cells.Item(17, 10).Value = 0.4448052188201075
$ws.Cells.Item(17, 12).Value = 0.7852521312778151
$ws.Cells.Item(17, 13).Value = 0.6046955779804506

$ws.Cells.Item(18, 2).Value = 1.928258136884267
$ws.Cells.Item(18, 4).Value = 0.1950547357559174
$ws.Cells.Item(18, 5).Value = 0.2883320338367739
$ws.Cells.Item(18, 6).Value = 2.023199358223707
$ws.Cells.Item(18, 7).Value = 1.515244872048754
$ws.Cells.Item(18, 8).Value = 1.316749907589866
$ws.Cells.Item(18, 10).Value = 0.4398749735485552
$ws.Cells.Item(18, 12).Value = 0.7699820991926742
$ws.Cells.Item(18, 13).Value = 0.5955847461880097

$ws.Cells.Item(19, 2).Value = 1.920769203348073
$ws.Cells.Item(19, 4).Value = 0.1944420811254304
$ws.Cells.Item(19, 5).Value = 0.2873046487468471
$ws.Cells.Item(19, 6).Value = 2.023789843799875
$ws.Cells.Item(19, 7).Value = 1.514089307326827
$ws.Cells.Item(19, 8).Value = 1.316865653432643
$ws.Cells.Item(19, 10).Value = 0.4382054927710328
$ws.Cells.Item(19, 12).Value = 0.7648120030964662
$ws.Cells.Item(19, 13).Value = 0.5925012891261474

$ws.Cells.Item(20, 2).Value = 1.954499665656272
$ws.Cells.Item(20, 4).Value = 0.1971968988860056
$ws.Cells.Item(20, 5).Value = 0.2919266006369199
$ws.Cells.Item(20, 6).Value = 2.021249352380025
$ws.Cells.Item(20, 7).Value = 1.519391415409416
$ws.Cells.Item(20, 8).Value = 1.316412788484342
$ws.Cells.Item(20, 10).Value = 0.4457176069472411
$ws.Cells.Item(20, 12).Value = 0.7880783038186223
$ws.Cells.Item(20, 13).Value = 0.6063824108632616

$ws.Cells.Item(21, 2).Value = 2.06834193249739
$ws.Cells.Item(21, 4).Value = 0.2064144983322507
$ws.Cells.Item(21, 5).Value = 0.3074308330082758
$ws.Cells.Item(21, 6).Value = 2.014763738074208
$ws.Cells.Item(21, 7).Value = 1.538998864386315
$ws.Cells.Item(21, 8).Value = 1.31608490435346
$ws.Cells.Item(21, 10).Value = 0.4709432910044882
$ws.Cells.Item(21, 12).Value = 0.8662534016417851
$ws.Cells.Item(21, 13).Value = 0.6531115177497071

$ws.Cells.Item(22, 2).Value = 2.143123549989582
$ws.Cells.Item(22, 4).Value = 0.2124105315266718
$ws.Cells.Item(22, 5).Value = 0.3175446327289393
$ws.Cells.Item(22, 6).Value = 2.012061224485421
$ws.Cells.Item(22, 7).Value = 1.55316101820091
$ws.Cells.Item(22, 8).Value = 1.316764660020141
$ws.Cells.Item(22, 10).Value = 0.4874179705275878
$ws.Cells.Item(22, 12).Value = 0.9173460849046364
$ws.Cells.Item(22, 13).Value = 0.6837179939571314

$ws.Cells.Item(23, 2).Value = 2.103176684198274
$ws.Cells.Item(23, 4).Value = 0.2092128978934369
$ws.Cells.Item(23, 5).Value = 0.3121484913677648
$ws.Cells.Item(23, 6).Value = 2.013362512029147
$ws.Cells.Item(23, 7).Value = 1.545478391770587
$ws.Cells.Item(23, 8).Value = 1.316319773770658
$ws.Cells.Item(23, 10).Value = 0.4786262845991587
$ws.Cells.Item(23, 12).Value = 0.8900770077092375
$ws.Cells.Item(23, 13).Value = 0.6673767433692603

$ws.Cells.Item(24, 2).Value = 1.952645414094661
$ws.Cells.Item(24, 4).Value = 0.1970457618989059
$ws.Cells.Item(24, 5).Value = 0.2916728773103543
$ws.Cells.Item(24, 6).Value = 2.021381172500213
$ws.Cells.Item(24, 7).Value = 1.519093531593683
$ws.Cells.Item(24, 8).Value = 1.316433178122111
$ws.Cells.Item(24, 10).Value = 0.4453051266937962
$ws.Cells.Item(24, 12).Value = 0.7868006115675144
$ws.Cells.Item(24, 13).Value = 0.605619782844002

$ws.Cells.Item(25, 2).Value = 1.791919210904098
$ws.Cells.Item(25, 4).Value = 0.1837986377970964
$ws.Cells.Item(25, 5).Value = 0.2695081643191699
$ws.Cells.Item(25, 6).Value = 2.036572663194136
$ws.Cells.Item(25, 7).Value = 1.496346481088267
$ws.Cells.Item(25, 8).Value = 1.320365599255496
$ws.Cells.Item(25, 10).Value = 0.4093210053905239
$ws.Cells.Item(25, 12).Value = 0.6754155471751346
$ws.Cells.Item(25, 13).Value = 0.5392966291227665

Write-Output "Done updating pl_mw sheet values."